$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Header row (row 1): new columns H/I/J ---
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Match the bold / bordered / centered style already used by the other
# header cells (B1:G1) without disturbing the existing shared style index.
$ws.Range("H1:J1").Font.Bold = $true
$ws.Range("H1:J1").HorizontalAlignment = -4108
$ws.Range("H1:J1").VerticalAlignment = -4160
$ws.Range("H1:J1").Borders.LineStyle = 1

# --- Data rows (rows 2-7): legislator_id / legislator_name / date values ---
$rows = 2,3,4,5,6,7
foreach ($r in $rows) {
    $hCell = $ws.Range("H$r")
    # Force text so Excel doesn't reinterpret the string as a date serial.
    $hCell.NumberFormat = "@"
    $hCell.Value = "2011-11-22"
    $hCell.ClearFormats()

    $ws.Range("I$r").Value = "費鴻泰"
    $ws.Range("J$r").Value = 1365
}
